$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Determine used range bounds
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

# Column G is "Recorded By"
$col = 7

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value()
    if ($null -eq $val) { continue }
    $text = [string]$val
    if ($text -eq "") { continue }

    $parts = $text -split ", "
    $n = $parts.Length
    if ($n -gt 1 -and $parts[$n - 1] -eq "System") {
        $reversedParts = @()
        for ($i = $n - 1; $i -ge 0; $i--) {
            $reversedParts += $parts[$i]
        }
        $newText = [string]::Join(", ", $reversedParts)
        $cell.Value = $newText
    }
}
